$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 114, shifting rows 114-119 down to 115-120.
$ws.Rows("114:114").Insert()

# Fill in the new row 114 with the new weekly data entry.
$ws.Cells.Item(114, 1).Value = 6
$ws.Cells.Item(114, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(114, 3).Value = "Metropolitana"
$ws.Cells.Item(114, 4).Value = 44509
$ws.Cells.Item(114, 4).NumberFormat = $ws.Cells.Item(115, 4).NumberFormat
$ws.Cells.Item(114, 5).Value = 13
$ws.Cells.Item(114, 6).Value = 100112001
$ws.Cells.Item(114, 7).Value = "Berenjena"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 180
$ws.Cells.Item(114, 11).Value = 11000
$ws.Cells.Item(114, 12).Value = 12000
$ws.Cells.Item(114, 13).Value = 11556
$ws.Cells.Item(114, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(114, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(114, 16).Value = 165
$ws.Cells.Item(114, 17).Value = 70
$ws.Cells.Item(114, 18).Value = "Hortaliza"
